$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 199; $r++) {
    $cell = $ws.Cells($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}
